$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark numeric-looking price cells as Text so COM keeps them as strings (matches source inline-string cells)
$ws.Range("D5,D6,D9,D10,D12,D14,D19,D20,D21,D23,D24,D26,D27,D29,D32,D33,D35,D36,D37,D40,D44,D45,D46,D48,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.106.69"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.952.66"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "380.33"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "102.14"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "36.50"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "0.0840"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "3.420.96"
$ws.Range("D14").Value = "18.03"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "2.941.98"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "51.065.35"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "3.19"
$ws.Range("E19").Value = "  -6.13%  "
$ws.Range("D20").Value = "7.10"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "68.45"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "261.52"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "8.34"
$ws.Range("E26").Value = "  +13.14%  "
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "4.10"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +8.88%  "
$ws.Range("D32").Value = "25.62"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "9.80"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "50.40"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "33.83"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "16.77"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "121.63"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "21.14"
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "2.003.84"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "3.22"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "0.0338"
$ws.Range("E51").Value = "  +4.85%  "
